# Append: 2025-10-10 18:31 JST
# Update the "取得日時" (acquired-at) timestamp in column A for every
# existing data row of the active sheet (ランサーズ) to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-10 18:31:23"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
